# Update "南宁-漫展信息" workbook per commit diff:
#  - Sheet "展览": F4 4860 -> 4876; insert a new row (广西·THO04-永夜廻想) before the
#    "花海演绎" row, shifting the remaining rows down, and bump F values on the
#    two shifted rows that also changed.
#  - Sheet "全部类型": same F4 bump, plus a new row for THO04-永夜廻想 inserted
#    before "花海演绎", with the same downstream shifts/bumps.

$wb = $excel.ActiveWorkbook

function Set-TextCell {
    param($ws, [string]$addr, [string]$text)
    # Force text storage (avoid Excel auto-detecting dates/numbers), then
    # strip the resulting "quote prefix" style artifact so formatting stays
    # identical to a plain unformatted text cell.
    $ws.Range($addr).Value = "'" + $text
    $ws.Range($addr).Style = "Normal"
}

# ---------------------------------------------------------------------
# Sheet 1: 展览
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F4").Value = 4876

# Insert a new row at position 5; existing rows 5-7 shift down to 6-8.
$ws1.Rows.Item(5).Insert()

# The inserted row's column-A cell needs the same number-style (bold,
# centered, bordered) as the rest of column A.
$srcA = $ws1.Cells.Item(4, 1)
$dstA = $ws1.Cells.Item(5, 1)
$dstA.Font.Bold = $srcA.Font.Bold
$dstA.HorizontalAlignment = $srcA.HorizontalAlignment
$dstA.VerticalAlignment = $srcA.VerticalAlignment
$dstA.Borders.LineStyle = $srcA.Borders.LineStyle

$ws1.Range("A5").Value = 4
Set-TextCell $ws1 "B5" "2024-10-04"
$ws1.Range("C5").Value = "广西·THO04-永夜廻想"
$ws1.Range("D5").Value = "北湖北路48-5号(近北湖小区) 金御华尊国际大酒店"
$ws1.Range("E5").Value = "2024.10.04 09:00-10.04 22:00"
$ws1.Range("F5").Value = 14
$ws1.Range("G5").Value = 60
$ws1.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=92574"
$ws1.Range("I5").Value = "//i2.hdslb.com/bfs/openplatform/202409/T5Qnv1zR1726732010464.jpeg"

# Row 6 (was row 5, 花海演绎): index + "想去人数" both changed.
$ws1.Range("A6").Value = 5
$ws1.Range("F6").Value = 21

# Row 7 (was row 6, 熊喵M): index + "想去人数" both changed.
$ws1.Range("A7").Value = 6
$ws1.Range("F7").Value = 47

# Row 8 (was row 7, 万圣漫控嘉年华10): only the index changes.
$ws1.Range("A8").Value = 7

# ---------------------------------------------------------------------
# Sheet 2: 全部类型
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F4").Value = 4876

# Insert a new row at position 6 (after 莫西干人, before 花海演绎); existing
# rows 6-9 shift down to 7-10.
$ws4.Rows.Item(6).Insert()

$srcA4 = $ws4.Cells.Item(5, 1)
$dstA4 = $ws4.Cells.Item(6, 1)
$dstA4.Font.Bold = $srcA4.Font.Bold
$dstA4.HorizontalAlignment = $srcA4.HorizontalAlignment
$dstA4.VerticalAlignment = $srcA4.VerticalAlignment
$dstA4.Borders.LineStyle = $srcA4.Borders.LineStyle

$ws4.Range("A6").Value = 5
Set-TextCell $ws4 "B6" "2024-10-04"
$ws4.Range("C6").Value = "广西·THO04-永夜廻想"
$ws4.Range("D6").Value = "北湖北路48-5号(近北湖小区) 金御华尊国际大酒店"
$ws4.Range("E6").Value = "2024.10.04 09:00-10.04 22:00"
$ws4.Range("F6").Value = 14
$ws4.Range("G6").Value = 60
$ws4.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=92574"
$ws4.Range("I6").Value = "//i2.hdslb.com/bfs/openplatform/202409/T5Qnv1zR1726732010464.jpeg"

# Row 7 (was row 6, 花海演绎): index + "想去人数" both changed.
$ws4.Range("A7").Value = 6
$ws4.Range("F7").Value = 21

# Row 8 (was row 7, 井草圣二): only the index changes.
$ws4.Range("A8").Value = 7

# Row 9 (was row 8, 熊喵M): index + "想去人数" both changed.
$ws4.Range("A9").Value = 8
$ws4.Range("F9").Value = 47

# Row 10 (was row 9, 万圣漫控嘉年华10): only the index changes.
$ws4.Range("A10").Value = 9
